$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.824.46'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +5.78%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.536.08'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +9.30%  '

$ws.Range("E4").Value = '  +0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '186.81'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +8.67%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '550.66'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.36%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.531.91'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +9.12%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.607'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.18%  '

$ws.Range("E9").Value = '  -0.01%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.629'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.84%  '

$ws.Range("E11").Value = '  +13.41%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '54.33'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.15%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000268'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +5.59%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.32'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.10%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.125.95'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +9.96%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.552.98'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +9.92%  '

$ws.Range("E17").Value = '  +4.12%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '66.966.15'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +6.21%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '18.12'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +5.34%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.89'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +7.52%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.990'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.26%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '430.08'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +17.38%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.04'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +7.05%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '84.78'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.35%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.08'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.72%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.02'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.43%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.88'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +9.08%  '

$ws.Range("E28").Value = '  -0.35%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '12.07'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +6.86%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.03'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +10.33%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '30.19'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.99%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '645.46'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.23%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.55'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.29%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '11.70'
$ws.Range("D34").Style = "Normal"

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.111'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.45%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '59.47'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.41%  '

$ws.Range("E37").Value = '  +22.68%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0₃0817'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +15.14%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '38.35'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.43%  '

$ws.Range("E40").Value = '  -0.20%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.389'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.02%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.36'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +14.58%  '

$ws.Range("E43").Value = '  +0.17%  '

$ws.Range("B44").Value = 'Maker'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.031.64'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +5.05%  '

$ws.Range("B45").Value = 'Fetch.AI'
$ws.Range("C45").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.67'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.76%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.38'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +11.54%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.87'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +10.96%  '

$ws.Range("B48").Value = 'WEMIXToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.79'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.80%  '

$ws.Range("B49").Value = 'VeChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0417'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +6.00%  '

$ws.Range("B50").Value = 'Stellar'
$ws.Range("C50").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.130'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +4.42%  '

$ws.Range("B51").Value = 'Monero'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '141.60'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +5.24%  '

